$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Swap row pairs (each pair exchanges B,F,G,H,I,J,K:AC) ----
# pair 47 <-> 48
$ws.Range("B47").Value2 = 6830649
$ws.Range("F47").Value2 = 'Toluca Women'
$ws.Range("G47").Value2 = 'Juarez FC Women'
$ws.Range("H47").Value2 = 1
$ws.Range("I47").Value2 = 1
$ws.Range("J47").Value2 = 'D'
$ws.Range("K47").Value2 = 2.3
$ws.Range("L47").Value2 = 3.2
$ws.Range("M47").Value2 = 2.7
$ws.Range("N47").Value2 = 2.15
$ws.Range("O47").Value2 = 3.2
$ws.Range("P47").Value2 = 3
$ws.Range("Q47").Value2 = -0.25
$ws.Range("R47").Value2 = 1.925
$ws.Range("S47").Value2 = 1.875
$ws.Range("T47").Value2 = 2.75
$ws.Range("U47").Value2 = 1.85
$ws.Range("V47").Value2 = 1.95
$ws.Range("W47").Value2 = -1
$ws.Range("X47").Value2 = 2.2
$ws.Range("Y47").Value2 = -1
$ws.Range("Z47").Value2 = -0.5
$ws.Range("AA47").Value2 = 0.4375
$ws.Range("AB47").Value2 = -1
$ws.Range("AC47").Value2 = 0.95
$ws.Range("B48").Value2 = 6830648
$ws.Range("F48").Value2 = 'Puebla Women'
$ws.Range("G48").Value2 = 'Unam Pumas Women'
$ws.Range("H48").Value2 = 1
$ws.Range("I48").Value2 = 1
$ws.Range("J48").Value2 = 'D'
$ws.Range("K48").Value2 = 6
$ws.Range("L48").Value2 = 4.333
$ws.Range("M48").Value2 = 1.4
$ws.Range("N48").Value2 = 7
$ws.Range("O48").Value2 = 4.75
$ws.Range("P48").Value2 = 1.333
$ws.Range("Q48").Value2 = 1.5
$ws.Range("R48").Value2 = 1.825
$ws.Range("S48").Value2 = 1.975
$ws.Range("T48").Value2 = 3
$ws.Range("U48").Value2 = 1.925
$ws.Range("V48").Value2 = 1.875
$ws.Range("W48").Value2 = -1
$ws.Range("X48").Value2 = 3.75
$ws.Range("Y48").Value2 = -1
$ws.Range("Z48").Value2 = 0.825
$ws.Range("AA48").Value2 = -1
$ws.Range("AB48").Value2 = -1
$ws.Range("AC48").Value2 = 0.875

# pair 101 <-> 102
$ws.Range("B101").Value2 = 6830703
$ws.Range("F101").Value2 = 'Leon Women'
$ws.Range("G101").Value2 = 'Atletico San Luis Women'
$ws.Range("H101").Value2 = 3
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 'H'
$ws.Range("K101").Value2 = 1.8
$ws.Range("L101").Value2 = 3.4
$ws.Range("M101").Value2 = 3.75
$ws.Range("N101").Value2 = 1.615
$ws.Range("O101").Value2 = 3.75
$ws.Range("P101").Value2 = 4.5
$ws.Range("Q101").Value2 = -1
$ws.Range("R101").Value2 = 2.025
$ws.Range("S101").Value2 = 1.775
$ws.Range("T101").Value2 = 3.25
$ws.Range("U101").Value2 = 2.025
$ws.Range("V101").Value2 = 1.775
$ws.Range("W101").Value2 = 0.615
$ws.Range("X101").Value2 = -1
$ws.Range("Y101").Value2 = -1
$ws.Range("Z101").Value2 = 1.025
$ws.Range("AA101").Value2 = -1
$ws.Range("AB101").Value2 = -0.5
$ws.Range("AC101").Value2 = 0.3875
$ws.Range("B102").Value2 = 7163155
$ws.Range("F102").Value2 = 'Tijuana Women'
$ws.Range("G102").Value2 = 'Atlas Women'
$ws.Range("H102").Value2 = 1
$ws.Range("I102").Value2 = 1
$ws.Range("J102").Value2 = 'D'
$ws.Range("K102").Value2 = 1.571
$ws.Range("L102").Value2 = 4
$ws.Range("M102").Value2 = 4.333
$ws.Range("N102").Value2 = 1.5
$ws.Range("O102").Value2 = 4.2
$ws.Range("P102").Value2 = 5
$ws.Range("Q102").Value2 = -1.25
$ws.Range("R102").Value2 = 1.975
$ws.Range("S102").Value2 = 1.825
$ws.Range("T102").Value2 = 3.25
$ws.Range("U102").Value2 = 1.85
$ws.Range("V102").Value2 = 1.95
$ws.Range("W102").Value2 = -1
$ws.Range("X102").Value2 = 3.2
$ws.Range("Y102").Value2 = -1
$ws.Range("Z102").Value2 = -1
$ws.Range("AA102").Value2 = 0.825
$ws.Range("AB102").Value2 = -1
$ws.Range("AC102").Value2 = 0.95

# pair 109 <-> 110
$ws.Range("B109").Value2 = 6830711
$ws.Range("F109").Value2 = 'Atletico San Luis Women'
$ws.Range("G109").Value2 = 'Tijuana Women'
$ws.Range("H109").Value2 = 0
$ws.Range("I109").Value2 = 2
$ws.Range("J109").Value2 = 'A'
$ws.Range("K109").Value2 = 3
$ws.Range("L109").Value2 = 3.6
$ws.Range("M109").Value2 = 2
$ws.Range("N109").Value2 = 4
$ws.Range("O109").Value2 = 3.8
$ws.Range("P109").Value2 = 1.666
$ws.Range("Q109").Value2 = 0.75
$ws.Range("R109").Value2 = 1.925
$ws.Range("S109").Value2 = 1.875
$ws.Range("T109").Value2 = 3
$ws.Range("U109").Value2 = 1.8
$ws.Range("V109").Value2 = 2
$ws.Range("W109").Value2 = -1
$ws.Range("X109").Value2 = -1
$ws.Range("Y109").Value2 = 0.6659999999999999
$ws.Range("Z109").Value2 = -1
$ws.Range("AA109").Value2 = 0.875
$ws.Range("AB109").Value2 = -1
$ws.Range("AC109").Value2 = 1
$ws.Range("B110").Value2 = 6830712
$ws.Range("F110").Value2 = 'Tigres UANL Women'
$ws.Range("G110").Value2 = 'Unam Pumas Women'
$ws.Range("H110").Value2 = 3
$ws.Range("I110").Value2 = 0
$ws.Range("J110").Value2 = 'H'
$ws.Range("K110").Value2 = 1.181
$ws.Range("L110").Value2 = 6
$ws.Range("M110").Value2 = 10
$ws.Range("N110").Value2 = 1.1
$ws.Range("O110").Value2 = 9
$ws.Range("P110").Value2 = 19
$ws.Range("Q110").Value2 = -2.5
$ws.Range("R110").Value2 = 1.8
$ws.Range("S110").Value2 = 2
$ws.Range("T110").Value2 = 3.75
$ws.Range("U110").Value2 = 1.8
$ws.Range("V110").Value2 = 2
$ws.Range("W110").Value2 = 0.1000000000000001
$ws.Range("X110").Value2 = -1
$ws.Range("Y110").Value2 = -1
$ws.Range("Z110").Value2 = 0.8
$ws.Range("AA110").Value2 = -1
$ws.Range("AB110").Value2 = -1
$ws.Range("AC110").Value2 = 1

# pair 131 <-> 132
$ws.Range("B131").Value2 = 6830732
$ws.Range("F131").Value2 = 'Cruz Azul Women'
$ws.Range("G131").Value2 = 'Tijuana Women'
$ws.Range("H131").Value2 = 3
$ws.Range("I131").Value2 = 6
$ws.Range("J131").Value2 = 'A'
$ws.Range("K131").Value2 = 3.2
$ws.Range("L131").Value2 = 3.5
$ws.Range("M131").Value2 = 2
$ws.Range("N131").Value2 = 3.75
$ws.Range("O131").Value2 = 3.5
$ws.Range("P131").Value2 = 1.833
$ws.Range("Q131").Value2 = 0.5
$ws.Range("R131").Value2 = 1.95
$ws.Range("S131").Value2 = 1.85
$ws.Range("T131").Value2 = 2.75
$ws.Range("U131").Value2 = 1.95
$ws.Range("V131").Value2 = 1.85
$ws.Range("W131").Value2 = -1
$ws.Range("X131").Value2 = -1
$ws.Range("Y131").Value2 = 0.833
$ws.Range("Z131").Value2 = -1
$ws.Range("AA131").Value2 = 0.8500000000000001
$ws.Range("AB131").Value2 = 0.95
$ws.Range("AC131").Value2 = -1
$ws.Range("B132").Value2 = 6830733
$ws.Range("F132").Value2 = 'Puebla Women'
$ws.Range("G132").Value2 = 'Club America Women'
$ws.Range("H132").Value2 = 1
$ws.Range("I132").Value2 = 6
$ws.Range("J132").Value2 = 'A'
$ws.Range("K132").Value2 = 29
$ws.Range("L132").Value2 = 11
$ws.Range("M132").Value2 = 1.062
$ws.Range("N132").Value2 = 34
$ws.Range("O132").Value2 = 13
$ws.Range("P132").Value2 = 1.045
$ws.Range("Q132").Value2 = 3.25
$ws.Range("R132").Value2 = 1.9
$ws.Range("S132").Value2 = 1.9
$ws.Range("T132").Value2 = 4
$ws.Range("U132").Value2 = 1.825
$ws.Range("V132").Value2 = 1.975
$ws.Range("W132").Value2 = -1
$ws.Range("X132").Value2 = -1
$ws.Range("Y132").Value2 = 0.04499999999999993
$ws.Range("Z132").Value2 = -1
$ws.Range("AA132").Value2 = 0.8999999999999999
$ws.Range("AB132").Value2 = 0.825
$ws.Range("AC132").Value2 = -1

# pair 149 <-> 150
$ws.Range("B149").Value2 = 6830751
$ws.Range("F149").Value2 = 'Leon Women'
$ws.Range("G149").Value2 = 'Atlas Women'
$ws.Range("H149").Value2 = 1
$ws.Range("I149").Value2 = 1
$ws.Range("J149").Value2 = 'D'
$ws.Range("K149").Value2 = 2.1
$ws.Range("L149").Value2 = 3.5
$ws.Range("M149").Value2 = 2.875
$ws.Range("N149").Value2 = 1.7
$ws.Range("O149").Value2 = 3.75
$ws.Range("P149").Value2 = 4
$ws.Range("Q149").Value2 = -0.75
$ws.Range("R149").Value2 = 1.95
$ws.Range("S149").Value2 = 1.85
$ws.Range("T149").Value2 = 3
$ws.Range("U149").Value2 = 1.925
$ws.Range("V149").Value2 = 1.875
$ws.Range("W149").Value2 = -1
$ws.Range("X149").Value2 = 2.75
$ws.Range("Y149").Value2 = -1
$ws.Range("Z149").Value2 = -1
$ws.Range("AA149").Value2 = 0.8500000000000001
$ws.Range("AB149").Value2 = -1
$ws.Range("AC149").Value2 = 0.875
$ws.Range("B150").Value2 = 6830750
$ws.Range("F150").Value2 = 'Club Necaxa Women'
$ws.Range("G150").Value2 = 'Tigres UANL Women'
$ws.Range("H150").Value2 = 1
$ws.Range("I150").Value2 = 3
$ws.Range("J150").Value2 = 'A'
$ws.Range("K150").Value2 = 23
$ws.Range("L150").Value2 = 13
$ws.Range("M150").Value2 = 1.03
$ws.Range("N150").Value2 = 29
$ws.Range("O150").Value2 = 12
$ws.Range("P150").Value2 = 1.055
$ws.Range("Q150").Value2 = 3
$ws.Range("R150").Value2 = 2
$ws.Range("S150").Value2 = 1.8
$ws.Range("T150").Value2 = 4
$ws.Range("U150").Value2 = 1.85
$ws.Range("V150").Value2 = 1.95
$ws.Range("W150").Value2 = -1
$ws.Range("X150").Value2 = -1
$ws.Range("Y150").Value2 = 0.05499999999999994
$ws.Range("Z150").Value2 = 1
$ws.Range("AA150").Value2 = -1
$ws.Range("AB150").Value2 = 0
$ws.Range("AC150").Value2 = -0

# pair 229 <-> 230
$ws.Range("B229").Value2 = 7645712
$ws.Range("F229").Value2 = 'Chivas Guadalajara Women'
$ws.Range("G229").Value2 = 'Santos Laguna Women'
$ws.Range("H229").Value2 = 10
$ws.Range("I229").Value2 = 2
$ws.Range("J229").Value2 = 'H'
$ws.Range("K229").Value2 = 1.1
$ws.Range("L229").Value2 = 8
$ws.Range("M229").Value2 = 13
$ws.Range("N229").Value2 = 1.03
$ws.Range("O229").Value2 = 17
$ws.Range("P229").Value2 = 41
$ws.Range("Q229").Value2 = -3.75
$ws.Range("R229").Value2 = 1.775
$ws.Range("S229").Value2 = 1.925
$ws.Range("T229").Value2 = 4.75
$ws.Range("U229").Value2 = 1.9
$ws.Range("V229").Value2 = 1.9
$ws.Range("W229").Value2 = 0.03000000000000003
$ws.Range("X229").Value2 = -1
$ws.Range("Y229").Value2 = -1
$ws.Range("Z229").Value2 = 0.7749999999999999
$ws.Range("AA229").Value2 = -1
$ws.Range("AB229").Value2 = 0.8999999999999999
$ws.Range("AC229").Value2 = -1
$ws.Range("B230").Value2 = 7645713
$ws.Range("F230").Value2 = 'Toluca Women'
$ws.Range("G230").Value2 = 'Tigres UANL Women'
$ws.Range("H230").Value2 = 0
$ws.Range("I230").Value2 = 7
$ws.Range("J230").Value2 = 'A'
$ws.Range("K230").Value2 = 9
$ws.Range("L230").Value2 = 7
$ws.Range("M230").Value2 = 1.166
$ws.Range("N230").Value2 = 5.25
$ws.Range("O230").Value2 = 4.75
$ws.Range("P230").Value2 = 1.4
$ws.Range("Q230").Value2 = 1.25
$ws.Range("R230").Value2 = 1.925
$ws.Range("S230").Value2 = 1.875
$ws.Range("T230").Value2 = 3
$ws.Range("U230").Value2 = 1.75
$ws.Range("V230").Value2 = 1.95
$ws.Range("W230").Value2 = -1
$ws.Range("X230").Value2 = -1
$ws.Range("Y230").Value2 = 0.3999999999999999
$ws.Range("Z230").Value2 = -1
$ws.Range("AA230").Value2 = 0.875
$ws.Range("AB230").Value2 = 0.75
$ws.Range("AC230").Value2 = -1

# pair 245 <-> 246
$ws.Range("B245").Value2 = 7645794
$ws.Range("F245").Value2 = 'Toluca Women'
$ws.Range("G245").Value2 = 'Club America Women'
$ws.Range("H245").Value2 = 3
$ws.Range("I245").Value2 = 0
$ws.Range("J245").Value2 = 'H'
$ws.Range("K245").Value2 = 6
$ws.Range("L245").Value2 = 5
$ws.Range("M245").Value2 = 1.333
$ws.Range("N245").Value2 = 5.75
$ws.Range("O245").Value2 = 5
$ws.Range("P245").Value2 = 1.363
$ws.Range("Q245").Value2 = 1.5
$ws.Range("R245").Value2 = 1.825
$ws.Range("S245").Value2 = 1.975
$ws.Range("T245").Value2 = 3.25
$ws.Range("U245").Value2 = 1.825
$ws.Range("V245").Value2 = 1.975
$ws.Range("W245").Value2 = 4.75
$ws.Range("X245").Value2 = -1
$ws.Range("Y245").Value2 = -1
$ws.Range("Z245").Value2 = 0.825
$ws.Range("AA245").Value2 = -1
$ws.Range("AB245").Value2 = -0.5
$ws.Range("AC245").Value2 = 0.4875
$ws.Range("B246").Value2 = 7645793
$ws.Range("F246").Value2 = 'Queretaro Women'
$ws.Range("G246").Value2 = 'Cruz Azul Women'
$ws.Range("H246").Value2 = 3
$ws.Range("I246").Value2 = 0
$ws.Range("J246").Value2 = 'H'
$ws.Range("K246").Value2 = 1.909
$ws.Range("L246").Value2 = 3.5
$ws.Range("M246").Value2 = 3.3
$ws.Range("N246").Value2 = 2.1
$ws.Range("O246").Value2 = 3.4
$ws.Range("P246").Value2 = 2.875
$ws.Range("Q246").Value2 = -0.25
$ws.Range("R246").Value2 = 1.875
$ws.Range("S246").Value2 = 1.925
$ws.Range("T246").Value2 = 2.75
$ws.Range("U246").Value2 = 2
$ws.Range("V246").Value2 = 1.8
$ws.Range("W246").Value2 = 1.1
$ws.Range("X246").Value2 = -1
$ws.Range("Y246").Value2 = -1
$ws.Range("Z246").Value2 = 0.875
$ws.Range("AA246").Value2 = -1
$ws.Range("AB246").Value2 = 0.5
$ws.Range("AC246").Value2 = -0.5

# pair 248 <-> 249
$ws.Range("B248").Value2 = 7645796
$ws.Range("F248").Value2 = 'Chivas Guadalajara Women'
$ws.Range("G248").Value2 = 'Club Necaxa Women'
$ws.Range("H248").Value2 = 4
$ws.Range("I248").Value2 = 0
$ws.Range("J248").Value2 = 'H'
$ws.Range("K248").Value2 = 1.055
$ws.Range("L248").Value2 = 10
$ws.Range("M248").Value2 = 21
$ws.Range("N248").Value2 = 1.062
$ws.Range("O248").Value2 = 11
$ws.Range("P248").Value2 = 29
$ws.Range("Q248").Value2 = -3
$ws.Range("R248").Value2 = 1.85
$ws.Range("S248").Value2 = 1.95
$ws.Range("T248").Value2 = 4
$ws.Range("U248").Value2 = 1.8
$ws.Range("V248").Value2 = 2
$ws.Range("W248").Value2 = 0.06200000000000006
$ws.Range("X248").Value2 = -1
$ws.Range("Y248").Value2 = -1
$ws.Range("Z248").Value2 = 0.8500000000000001
$ws.Range("AA248").Value2 = -1
$ws.Range("AB248").Value2 = 0
$ws.Range("AC248").Value2 = -0
$ws.Range("B249").Value2 = 7645719
$ws.Range("F249").Value2 = 'Monterrey Women'
$ws.Range("G249").Value2 = 'Santos Laguna Women'
$ws.Range("H249").Value2 = 6
$ws.Range("I249").Value2 = 0
$ws.Range("J249").Value2 = 'H'
$ws.Range("K249").Value2 = 1.025
$ws.Range("L249").Value2 = 15
$ws.Range("M249").Value2 = 34
$ws.Range("N249").Value2 = 1.01
$ws.Range("O249").Value2 = 34
$ws.Range("P249").Value2 = 67
$ws.Range("Q249").Value2 = -4.75
$ws.Range("R249").Value2 = 1.775
$ws.Range("S249").Value2 = 1.925
$ws.Range("T249").Value2 = 5.75
$ws.Range("U249").Value2 = 1.85
$ws.Range("V249").Value2 = 1.95
$ws.Range("W249").Value2 = 0.01000000000000001
$ws.Range("X249").Value2 = -1
$ws.Range("Y249").Value2 = -1
$ws.Range("Z249").Value2 = 0.7749999999999999
$ws.Range("AA249").Value2 = -1
$ws.Range("AB249").Value2 = 0.425
$ws.Range("AC249").Value2 = -0.5

# pair 263 <-> 264
$ws.Range("B263").Value2 = 7645807
$ws.Range("F263").Value2 = 'Club Necaxa Women'
$ws.Range("G263").Value2 = 'Leon Women'
$ws.Range("H263").Value2 = 2
$ws.Range("I263").Value2 = 1
$ws.Range("J263").Value2 = 'H'
$ws.Range("K263").Value2 = 4.333
$ws.Range("L263").Value2 = 4
$ws.Range("M263").Value2 = 1.571
$ws.Range("N263").Value2 = 7
$ws.Range("O263").Value2 = 4.2
$ws.Range("P263").Value2 = 1.363
$ws.Range("Q263").Value2 = 1.5
$ws.Range("R263").Value2 = 1.75
$ws.Range("S263").Value2 = 1.95
$ws.Range("T263").Value2 = 2.75
$ws.Range("U263").Value2 = 1.8
$ws.Range("V263").Value2 = 2
$ws.Range("W263").Value2 = 6
$ws.Range("X263").Value2 = -1
$ws.Range("Y263").Value2 = -1
$ws.Range("Z263").Value2 = 0.75
$ws.Range("AA263").Value2 = -1
$ws.Range("AB263").Value2 = 0.4
$ws.Range("AC263").Value2 = -0.5
$ws.Range("B264").Value2 = 7645806
$ws.Range("F264").Value2 = 'Atletico San Luis Women'
$ws.Range("G264").Value2 = 'Atlas Women'
$ws.Range("H264").Value2 = 0
$ws.Range("I264").Value2 = 0
$ws.Range("J264").Value2 = 'D'
$ws.Range("K264").Value2 = 3.2
$ws.Range("L264").Value2 = 3.6
$ws.Range("M264").Value2 = 1.909
$ws.Range("N264").Value2 = 3.1
$ws.Range("O264").Value2 = 3.6
$ws.Range("P264").Value2 = 2
$ws.Range("Q264").Value2 = 0.25
$ws.Range("R264").Value2 = 1.975
$ws.Range("S264").Value2 = 1.825
$ws.Range("T264").Value2 = 3
$ws.Range("U264").Value2 = 1.9
$ws.Range("V264").Value2 = 1.9
$ws.Range("W264").Value2 = -1
$ws.Range("X264").Value2 = 2.6
$ws.Range("Y264").Value2 = -1
$ws.Range("Z264").Value2 = 0.4875
$ws.Range("AA264").Value2 = -0.5
$ws.Range("AB264").Value2 = -1
$ws.Range("AC264").Value2 = 0.8999999999999999

# pair 271 <-> 272
$ws.Range("B271").Value2 = 7645809
$ws.Range("F271").Value2 = 'Mazatlan FC Women'
$ws.Range("G271").Value2 = 'Queretaro Women'
$ws.Range("H271").Value2 = 2
$ws.Range("I271").Value2 = 2
$ws.Range("J271").Value2 = 'D'
$ws.Range("K271").Value2 = 4.75
$ws.Range("L271").Value2 = 4
$ws.Range("M271").Value2 = 1.533
$ws.Range("N271").Value2 = 4.75
$ws.Range("O271").Value2 = 3.75
$ws.Range("P271").Value2 = 1.571
$ws.Range("Q271").Value2 = 1
$ws.Range("R271").Value2 = 1.775
$ws.Range("S271").Value2 = 2.025
$ws.Range("T271").Value2 = 2.75
$ws.Range("U271").Value2 = 1.95
$ws.Range("V271").Value2 = 1.85
$ws.Range("W271").Value2 = -1
$ws.Range("X271").Value2 = 2.75
$ws.Range("Y271").Value2 = -1
$ws.Range("Z271").Value2 = 0.7749999999999999
$ws.Range("AA271").Value2 = -1
$ws.Range("AB271").Value2 = 0.95
$ws.Range("AC271").Value2 = -1
$ws.Range("B272").Value2 = 7645812
$ws.Range("F272").Value2 = 'Atlas Women'
$ws.Range("G272").Value2 = 'Monterrey Women'
$ws.Range("H272").Value2 = 0
$ws.Range("I272").Value2 = 1
$ws.Range("J272").Value2 = 'A'
$ws.Range("K272").Value2 = 6.5
$ws.Range("L272").Value2 = 5
$ws.Range("M272").Value2 = 1.333
$ws.Range("N272").Value2 = 6.5
$ws.Range("O272").Value2 = 5
$ws.Range("P272").Value2 = 1.3
$ws.Range("Q272").Value2 = 1.5
$ws.Range("R272").Value2 = 1.875
$ws.Range("S272").Value2 = 1.925
$ws.Range("T272").Value2 = 3.25
$ws.Range("U272").Value2 = 2
$ws.Range("V272").Value2 = 1.8
$ws.Range("W272").Value2 = -1
$ws.Range("X272").Value2 = -1
$ws.Range("Y272").Value2 = 0.3
$ws.Range("Z272").Value2 = 0.875
$ws.Range("AA272").Value2 = -1
$ws.Range("AB272").Value2 = -1
$ws.Range("AC272").Value2 = 0.8
# ---- Append new rows 277-281 ----
# row 277
$ws.Range("A277").Value2 = 275
$ws.Range("B277").Value2 = 7645729
$ws.Range("C277").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("D277").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("E277").Value2 = 45380
$ws.Range("F277").Value2 = 'Monterrey Women'
$ws.Range("G277").Value2 = 'Tijuana Women'
$ws.Range("H277").Value2 = 3
$ws.Range("I277").Value2 = 1
$ws.Range("J277").Value2 = 'H'
$ws.Range("K277").Value2 = 1.285
$ws.Range("L277").Value2 = 5
$ws.Range("M277").Value2 = 7.5
$ws.Range("N277").Value2 = 1.2
$ws.Range("O277").Value2 = 7
$ws.Range("P277").Value2 = 8.5
$ws.Range("Q277").Value2 = -1.75
$ws.Range("R277").Value2 = 1.95
$ws.Range("S277").Value2 = 1.85
$ws.Range("T277").Value2 = 3.25
$ws.Range("U277").Value2 = 1.85
$ws.Range("V277").Value2 = 1.95
$ws.Range("W277").Value2 = 0.2
$ws.Range("X277").Value2 = -1
$ws.Range("Y277").Value2 = -1
$ws.Range("Z277").Value2 = 0.475
$ws.Range("AA277").Value2 = -0.5
$ws.Range("AB277").Value2 = 0.8500000000000001
$ws.Range("AC277").Value2 = -1

# row 278
$ws.Range("A278").Value2 = 276
$ws.Range("B278").Value2 = 7645815
$ws.Range("C278").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("D278").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("E278").Value2 = 45380.78125
$ws.Range("F278").Value2 = 'Cruz Azul Women'
$ws.Range("G278").Value2 = 'Club America Women'
$ws.Range("K278").Value2 = 10
$ws.Range("L278").Value2 = 8
$ws.Range("M278").Value2 = 1.125
$ws.Range("N278").Value2 = 19
$ws.Range("O278").Value2 = 9
$ws.Range("P278").Value2 = 1.1
$ws.Range("Q278").Value2 = 2.75
$ws.Range("R278").Value2 = 1.8
$ws.Range("S278").Value2 = 2
$ws.Range("T278").Value2 = 3.75
$ws.Range("U278").Value2 = 1.825
$ws.Range("V278").Value2 = 1.975
$ws.Range("W278").Value2 = 0
$ws.Range("X278").Value2 = 0
$ws.Range("Y278").Value2 = 0
$ws.Range("Z278").Value2 = 0
$ws.Range("AA278").Value2 = 0

# row 279
$ws.Range("A279").Value2 = 277
$ws.Range("B279").Value2 = 7645816
$ws.Range("C279").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("D279").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("E279").Value2 = 45380.83333333334
$ws.Range("F279").Value2 = 'Queretaro Women'
$ws.Range("G279").Value2 = 'Club Necaxa Women'
$ws.Range("K279").Value2 = 1.333
$ws.Range("L279").Value2 = 4.75
$ws.Range("M279").Value2 = 6.5
$ws.Range("N279").Value2 = 1.285
$ws.Range("O279").Value2 = 5.75
$ws.Range("P279").Value2 = 7.5
$ws.Range("Q279").Value2 = -1.75
$ws.Range("R279").Value2 = 2
$ws.Range("S279").Value2 = 1.8
$ws.Range("T279").Value2 = 3
$ws.Range("U279").Value2 = 1.825
$ws.Range("V279").Value2 = 1.975
$ws.Range("W279").Value2 = 0
$ws.Range("X279").Value2 = 0
$ws.Range("Y279").Value2 = 0
$ws.Range("Z279").Value2 = 0
$ws.Range("AA279").Value2 = 0

# row 280
$ws.Range("A280").Value2 = 278
$ws.Range("B280").Value2 = 7645817
$ws.Range("C280").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("D280").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("E280").Value2 = 45380.92083333333
$ws.Range("F280").Value2 = 'Leon Women'
$ws.Range("G280").Value2 = 'Mazatlan FC Women'
$ws.Range("K280").Value2 = 1.333
$ws.Range("L280").Value2 = 4.75
$ws.Range("M280").Value2 = 6.5
$ws.Range("N280").Value2 = 1.2
$ws.Range("O280").Value2 = 6
$ws.Range("P280").Value2 = 11
$ws.Range("Q280").Value2 = -2
$ws.Range("R280").Value2 = 1.9
$ws.Range("S280").Value2 = 1.9
$ws.Range("T280").Value2 = 3.5
$ws.Range("U280").Value2 = 1.95
$ws.Range("V280").Value2 = 1.85
$ws.Range("W280").Value2 = 0
$ws.Range("X280").Value2 = 0
$ws.Range("Y280").Value2 = 0
$ws.Range("Z280").Value2 = 0
$ws.Range("AA280").Value2 = 0

# row 281
$ws.Range("A281").Value2 = 279
$ws.Range("B281").Value2 = 7645820
$ws.Range("C281").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("D281").Value2 = 'Mexico Liga MX Femenil'
$ws.Range("E281").Value2 = 45381.83333333334
$ws.Range("F281").Value2 = 'Atletico San Luis Women'
$ws.Range("G281").Value2 = 'Pachuca Women'
$ws.Range("K281").Value2 = 9
$ws.Range("L281").Value2 = 4
$ws.Range("M281").Value2 = 1.333
$ws.Range("N281").Value2 = 8
$ws.Range("O281").Value2 = 4.2
$ws.Range("P281").Value2 = 1.333
$ws.Range("Q281").Value2 = 1.25
$ws.Range("R281").Value2 = 2.05
$ws.Range("S281").Value2 = 1.75
$ws.Range("T281").Value2 = 3.25
$ws.Range("U281").Value2 = 2.025
$ws.Range("V281").Value2 = 1.775
$ws.Range("W281").Value2 = 0
$ws.Range("X281").Value2 = 0
$ws.Range("Y281").Value2 = 0
$ws.Range("Z281").Value2 = 0
$ws.Range("AA281").Value2 = 0

# ---- Apply formatting to new rows 277-281 (copy from row 276's A and E cells) ----
$ws.Range("A276").Copy() | Out-Null
$ws.Range("A277:A281").PasteSpecial(-4122) | Out-Null
$ws.Range("E276").Copy() | Out-Null
$ws.Range("E277:E281").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
